$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H92").Value = 977.2381
$ws.Range("I92").Value = 909.0526
$ws.Range("J92").Value = 1625
$ws.Range("K92").Value = 909.0526
$ws.Range("L92").Value = 1625
$ws.Range("M92").Value = 338.9474
$ws.Range("N92").Value = -4121

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H96").Value = 1110.1111
$ws.Range("I96").Value = 755.1818
$ws.Range("J96").Value = 1667.8572
$ws.Range("K96").Value = 2265.5454
$ws.Range("L96").Value = 5003.571599999999
$ws.Range("M96").Value = -892.5454
$ws.Range("N96").Value = -7749.571599999999

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H97").Value = 1370
$ws.Range("I97").Value = 0
$ws.Range("J97").Value = 1370
$ws.Range("K97").Value = 0
$ws.Range("L97").Value = 4110
$ws.Range("N97").Value = -5102

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H111").Value = 6028.5
$ws.Range("I111").Value = 2242.6667
$ws.Range("J111").Value = 8300
$ws.Range("K111").Value = 6728.000100000001
$ws.Range("L111").Value = 24900
$ws.Range("M111").Value = -3661.000100000001

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H129").Value = 828.62
$ws.Range("I129").Value = 800
$ws.Range("J129").Value = 828.9091
$ws.Range("K129").Value = 2400
$ws.Range("L129").Value = 2486.7273
$ws.Range("M129").Value = 2600
$ws.Range("N129").Value = -12486.7273

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 1732.6857
$ws.Range("I138").Value = 540.6957
$ws.Range("J138").Value = 4017.3333
$ws.Range("K138").Value = 1622.0871
$ws.Range("L138").Value = 12051.9999
$ws.Range("M138").Value = 3517.9129
$ws.Range("N138").Value = -22331.9999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H9").Value = 0
$ws.Range("I9").Value = 0
$ws.Range("J9").Value = 0
$ws.Range("K9").Value = 0
$ws.Range("L9").Value = 0
$ws.Range("N9").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H20").Value = 0
$ws.Range("I20").Value = 0
$ws.Range("J20").Value = 0
$ws.Range("K20").Value = 0
$ws.Range("L20").Value = 0
$ws.Range("N20").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H23").Value = 21200
$ws.Range("I23").Value = 25000
$ws.Range("J23").Value = 19933.334
$ws.Range("K23").Value = 25000
$ws.Range("L23").Value = 19933.334
$ws.Range("M23").Value = -24741
$ws.Range("N23").Value = -20451.334

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2820.673
$ws.Range("I32").Value = 2214.8408
$ws.Range("J32").Value = 6152.75
$ws.Range("K32").Value = 2214.8408
$ws.Range("L32").Value = 6152.75
$ws.Range("M32").Value = -1927.8408

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H37").Value = 20530
$ws.Range("I37").Value = 2000
$ws.Range("J37").Value = 29795
$ws.Range("K37").Value = 2000
$ws.Range("L37").Value = 29795
$ws.Range("M37").Value = -1727
$ws.Range("N37").Value = -30341

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H44").Value = 22800
$ws.Range("I44").Value = 0
$ws.Range("J44").Value = 22800
$ws.Range("K44").Value = 0
$ws.Range("L44").Value = 22800
$ws.Range("N44").Value = -23776

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H51").Value = 18000
$ws.Range("I51").Value = 0
$ws.Range("J51").Value = 18000
$ws.Range("K51").Value = 0
$ws.Range("L51").Value = 18000
$ws.Range("N51").Value = -19512

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H55").Value = 18505.666
$ws.Range("I55").Value = 4500.5
$ws.Range("J55").Value = 25508.25
$ws.Range("K55").Value = 4500.5
$ws.Range("L55").Value = 25508.25
$ws.Range("M55").Value = -4185.5
$ws.Range("N55").Value = -26138.25

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 1003253.3
$ws.Range("I61").Value = 1386858.5
$ws.Range("J61").Value = 5880
$ws.Range("K61").Value = 1386858.5
$ws.Range("L61").Value = 5880
$ws.Range("M61").Value = -1386646.5
$ws.Range("N61").Value = -6304

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 2208.8948
$ws.Range("I74").Value = 2180.5293
$ws.Range("J74").Value = 2450
$ws.Range("K74").Value = 2180.5293
$ws.Range("L74").Value = 2450
$ws.Range("M74").Value = -1306.5293

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 2208.8948
$ws.Range("I77").Value = 2180.5293
$ws.Range("J77").Value = 2450
$ws.Range("K77").Value = 10902.6465
$ws.Range("L77").Value = 12250
$ws.Range("M77").Value = -6534.646500000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H102").Value = 2229.4211
$ws.Range("I102").Value = 859.9231
$ws.Range("J102").Value = 5196.6665
$ws.Range("K102").Value = 859.9231
$ws.Range("L102").Value = 5196.6665
$ws.Range("M102").Value = 762.0769
$ws.Range("N102").Value = -8440.666499999999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value = 1529.0555
$ws.Range("I122").Value = 1494.5
$ws.Range("J122").Value = 1650
$ws.Range("K122").Value = 4483.5
$ws.Range("L122").Value = 4950
$ws.Range("M122").Value = -2033.5
$ws.Range("N122").Value = -9850

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 20295
$ws.Range("I132").Value = 2268.3076
$ws.Range("J132").Value = 35918.133
$ws.Range("K132").Value = 6804.9228
$ws.Range("L132").Value = 107754.399
$ws.Range("M132").Value = -4274.9228
$ws.Range("N132").Value = -112814.399

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 1003253.3
$ws.Range("I136").Value = 1386858.5
$ws.Range("J136").Value = 5880
$ws.Range("K136").Value = 4160575.5
$ws.Range("L136").Value = 17640
$ws.Range("M136").Value = -4158025.5
$ws.Range("N136").Value = -22740

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 4484.857
$ws.Range("I94").Value = 2423.818
$ws.Range("J94").Value = 6752
$ws.Range("K94").Value = 2423.818
$ws.Range("L94").Value = 6752
$ws.Range("M94").Value = -1972.818
$ws.Range("N94").Value = -7654

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 2389.7273
$ws.Range("I99").Value = 2130.7
$ws.Range("J99").Value = 4980
$ws.Range("K99").Value = 2130.7
$ws.Range("L99").Value = 4980
$ws.Range("M99").Value = -632.6999999999998

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 9680.286
$ws.Range("I134").Value = 9680.286
$ws.Range("J134").Value = 0
$ws.Range("K134").Value = 29040.858
$ws.Range("L134").Value = 0
$ws.Range("M134").Value = -26505.858
$ws.Range("N134").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 10328.405
$ws.Range("I31").Value = 12955.593
$ws.Range("J31").Value = 3235
$ws.Range("K31").Value = 12955.593
$ws.Range("L31").Value = 3235
$ws.Range("M31").Value = -12660.593

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 10328.405
$ws.Range("I34").Value = 12955.593
$ws.Range("J34").Value = 3235
$ws.Range("K34").Value = 12955.593
$ws.Range("L34").Value = 3235
$ws.Range("M34").Value = -12753.593

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 24060.545
$ws.Range("I58").Value = 1477.4667
$ws.Range("J58").Value = 72452.86
$ws.Range("K58").Value = 1477.4667
$ws.Range("L58").Value = 72452.86
$ws.Range("M58").Value = -1274.4667
$ws.Range("N58").Value = -72858.86

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H122").Value = 1200
$ws.Range("I122").Value = 1200
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 3600
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -1150
$ws.Range("N122").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 20736.607
$ws.Range("I132").Value = 22733.334
$ws.Range("J132").Value = 8756.25
$ws.Range("K132").Value = 68200.00199999999
$ws.Range("L132").Value = 26268.75
$ws.Range("M132").Value = -65670.00199999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 1230.9565
$ws.Range("I134").Value = 956
$ws.Range("J134").Value = 1658.6666
$ws.Range("K134").Value = 2868
$ws.Range("L134").Value = 4975.9998
$ws.Range("M134").Value = -333
$ws.Range("N134").Value = -10045.9998

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H136").Value = 24060.545
$ws.Range("I136").Value = 1477.4667
$ws.Range("J136").Value = 72452.86
$ws.Range("K136").Value = 4432.4001
$ws.Range("L136").Value = 217358.58
$ws.Range("M136").Value = -1882.4001
$ws.Range("N136").Value = -222458.58

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H57").Value = 6500
$ws.Range("I57").Value = 3000
$ws.Range("J57").Value = 10000
$ws.Range("K57").Value = 9000
$ws.Range("L57").Value = 30000
$ws.Range("M57").Value = -8441
$ws.Range("N57").Value = -31118

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H107").Value = 12746.625
$ws.Range("I107").Value = 50050
$ws.Range("J107").Value = 312.16666
$ws.Range("K107").Value = 150150
$ws.Range("L107").Value = 936.4999799999999
$ws.Range("M107").Value = -148230
$ws.Range("N107").Value = -4776.49998

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H122").Value = 639.6667
$ws.Range("I122").Value = 393.85715
$ws.Range("J122").Value = 1500
$ws.Range("K122").Value = 3544.71435
$ws.Range("L122").Value = 13500
$ws.Range("M122").Value = -1094.71435
$ws.Range("N122").Value = -18400

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 721.37
$ws.Range("I131").Value = 240
$ws.Range("J131").Value = 741.42706
$ws.Range("K131").Value = 720
$ws.Range("L131").Value = 2224.28118
$ws.Range("M131").Value = 4320
$ws.Range("N131").Value = -12304.28118

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 3509.9
$ws.Range("I113").Value = 2750
$ws.Range("J113").Value = 4016.5
$ws.Range("K113").Value = 2750
$ws.Range("L113").Value = 4016.5
$ws.Range("M113").Value = -580
$ws.Range("N113").Value = -8356.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 33872.65
$ws.Range("I132").Value = 4639.727
$ws.Range("J132").Value = 87466.336
$ws.Range("K132").Value = 13919.181
$ws.Range("L132").Value = 262399.008
$ws.Range("M132").Value = -11389.181
$ws.Range("N132").Value = -267459.008

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 3371.2856
$ws.Range("I22").Value = 5150.5
$ws.Range("J22").Value = 2659.6
$ws.Range("K22").Value = 5150.5
$ws.Range("L22").Value = 2659.6
$ws.Range("M22").Value = -4855.5
$ws.Range("N22").Value = -3249.6

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H27").Value = 3371.2856
$ws.Range("I27").Value = 5150.5
$ws.Range("J27").Value = 2659.6
$ws.Range("K27").Value = 5150.5
$ws.Range("L27").Value = 2659.6
$ws.Range("M27").Value = -5043.5
$ws.Range("N27").Value = -2873.6

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H68").Value = 0
$ws.Range("I68").Value = 0
$ws.Range("J68").Value = 0
$ws.Range("K68").Value = 0
$ws.Range("L68").Value = 0
$ws.Range("N68").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H71").Value = 0
$ws.Range("I71").Value = 0
$ws.Range("J71").Value = 0
$ws.Range("K71").Value = 0
$ws.Range("L71").Value = 0
$ws.Range("N71").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 6244.3335
$ws.Range("I96").Value = 2000
$ws.Range("J96").Value = 7093.2
$ws.Range("K96").Value = 2000
$ws.Range("L96").Value = 7093.2
$ws.Range("M96").Value = -627
$ws.Range("N96").Value = -9839.200000000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H98").Value = 50000
$ws.Range("I98").Value = 0
$ws.Range("J98").Value = 50000
$ws.Range("K98").Value = 0
$ws.Range("L98").Value = 50000
$ws.Range("N98").Value = -55990

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 1261.0834
$ws.Range("I100").Value = 821.6667
$ws.Range("J100").Value = 1700.5
$ws.Range("K100").Value = 1643.3334
$ws.Range("L100").Value = 3401
$ws.Range("M100").Value = -1102.3334

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 1647.6666
$ws.Range("I122").Value = 1479.5714
$ws.Range("J122").Value = 2236
$ws.Range("K122").Value = 4438.7142
$ws.Range("L122").Value = 6708
$ws.Range("M122").Value = -1988.7142
$ws.Range("N122").Value = -11608

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 1793.9
$ws.Range("I126").Value = 1034.75
$ws.Range("J126").Value = 2300
$ws.Range("K126").Value = 3104.25
$ws.Range("L126").Value = 6900
$ws.Range("M126").Value = -634.25
$ws.Range("N126").Value = -11840

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 3447.8333
$ws.Range("I132").Value = 2839.7144
$ws.Range("J132").Value = 4299.2
$ws.Range("K132").Value = 8519.143199999999
$ws.Range("L132").Value = 12897.6
$ws.Range("M132").Value = -5989.143199999999
$ws.Range("N132").Value = -17957.6
